# Rewrote preprocessing for better accounting of timings.
#
# "glycemie moyenne estimee" was an excluded-criterion entry mistakenly
# listed among the *included* variables (column A, row 8). Move it over to
# the "excluded" column (column B), next to the other excluded lab values,
# and close the gap it leaves in column A by shifting every following
# included-variable row up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$movedValue = "glycemie moyenne estimee"
$sourceRow = 8
$destRow = 10

# Figure out how far column A currently extends.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row() + $usedRange.Rows.Count() - 1

# Capture the column-A values below the row being removed, in order, so we
# can shift them up by one once the source row is gone.
$colAValues = @()
for ($r = $sourceRow + 1; $r -le $lastRow; $r++) {
    $colAValues += $ws.Cells.Item($r, 1).Value()
}

for ($i = 0; $i -lt $colAValues.Length; $i++) {
    $ws.Cells.Item($sourceRow + $i, 1).Value = $colAValues[$i]
}

# The old final row is now a duplicate of the row above it; clear it so the
# sheet's used range shrinks by one row.
$ws.Cells.Item($lastRow, 1).Value = $null

# Place the moved entry into column B, alongside the other excluded values.
$ws.Cells.Item($destRow, 2).Value = $movedValue

# Restore the view state (scroll position / active selection) seen after
# the edit.
$ws.Range("C68").Select()
$excel.ActiveWindow.ScrollRow = 61
